# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets:
#   F2: 529 -> 530
#   F4: 17  -> 18
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 530
    $ws.Range("F4").Value = 18
}
